$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 2.04
$ws.Range("H2").Value = 3.65
$ws.Range("N2").Value = 6.4
$ws.Range("P2").Value = 2.88
$ws.Range("R2").Value = 1.77
$ws.Range("S2").Value = 2.16
$ws.Range("U2").Value = 2.74
$ws.Range("W2").Value = 1.97
$ws.Range("AK2").Value = 18.5
$ws.Range("H3").Value = 1.57
$ws.Range("I3").Value = 1.68
$ws.Range("I5").Value = 2.66
$ws.Range("F7").Value = 1.96
$ws.Range("Q7").Value = 1.53
$ws.Range("G8").Value = 2.74
$ws.Range("Q8").Value = 1.95
$ws.Range("H9").Value = 1.49
$ws.Range("I9").Value = 1.6
$ws.Range("Q9").Value = 1.57
$ws.Range("G10").Value = 2.56
$ws.Range("Q10").Value = 1.84
$ws.Range("Q11").Value = 1.78
$ws.Range("G14").Value = 2.9
$ws.Range("P14").Value = 1.78
$ws.Range("F15").Value = 5.2
$ws.Range("J16").Value = 3.85
$ws.Range("H17").Value = 1.68
$ws.Range("J17").Value = 3.7
$ws.Range("G18").Value = 1.26
$ws.Range("H18").Value = 13.5
$ws.Range("I18").Value = 16.5
$ws.Range("J18").Value = 7.4
$ws.Range("P18").Value = 3.6
$ws.Range("R18").Value = 2.04
$ws.Range("T18").Value = 1.69
$ws.Range("U18").Value = 2.16
$ws.Range("AA18").Value = 640
$ws.Range("AE18").Value = 190
$ws.Range("AH18").Value = 30
$ws.Range("AJ18").Value = 14.5
$ws.Range("AN18").Value = 3.15
$ws.Range("AO18").Value = 170
$ws.Range("AO19").Value = 19
